$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(40, "2023-03-21-12:00", 28026.1, 28252.5, 27907, 28167.7, 27941.04896359393, 28119.24264361145, 45890.77)
    ,@(41, "2023-03-21-13:00", 28167.6, 28200, 27872.1, 28019.1, 28119.24264361145, 28163.20461803625, 45961.945)
    ,@(42, "2023-03-21-14:00", 28019.2, 28108.5, 27706, 27870.5, 28163.20461803625, 28020.013496295, 49027.562)
    ,@(43, "2023-03-21-15:00", 27870.5, 28091.6, 27838.2, 28028.7, 28020.013496295, 28100.71539929858, 25473.391)
    ,@(44, "2023-03-21-16:00", 28028.7, 28460, 27998, 28334.9, 28100.71539929858, 28273.1581835965, 64811.513)
    ,@(45, "2023-03-21-17:00", 28335, 28412.4, 27944, 27979.9, 28273.1581835965, 28145.89490835118, 56073.654)
    ,@(46, "2023-03-21-18:00", 27979.8, 28149, 27921.5, 28119.3, 28145.89490835118, 28326.42688768086, 29461.273)
    ,@(47, "2023-03-21-19:00", 28118.9, 28217, 28071.1, 28111.3, 28326.42688768086, 28178.66962040419, 15047.394)
    ,@(48, "2023-03-21-20:00", 28111.2, 28200, 28020, 28058, 28178.66962040419, 28112.19943414073, 14765.408)
    ,@(49, "2023-03-21-21:00", 28058, 28085, 27902, 27990, 28112.19943414073, 28094.83401953691, 15159.448)
    ,@(50, "2023-03-21-22:00", 27990, 28053, 27910, 28019.1, 28094.83401953691, 27967.14361037849, 10237.152)
    ,@(51, "2023-03-21-23:00", 28019.2, 28137, 27970, 28091.1, 27967.14361037849, 28029.75347543266, 10417.92)
    ,@(52, "2023-03-22-00:00", 28091.1, 28099.7, 27930, 28015, 28029.75347543266, 28088.50912627042, 14113.985)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
